$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$p55 = $d.Paragraphs.Item($count - 1)   # trailing "a3"-styled empty paragraph
$p56 = $d.Paragraphs.Item($count)       # trailing bare empty paragraph right before sectPr

# Step 1: update the existing trailing "a3" paragraph in place -
# drop the eastAsia font entry and mark the (now empty) run language as en-US.
$xmlP55 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a3"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$p55.Range.InsertXML($xmlP55)

# Step 2: add a new paragraph right after it holding the new note about the
# password-change bug that was found/fixed during this sprint.
$p55 = $d.Paragraphs.Item($count - 1)
$p55.Range.InsertParagraphAfter()
$textP = $d.Paragraphs.Item($count)
$xmlText = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a3"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Κατά την διάρκεια των αλλαγών της εφαρμογής ολοκληρώθηκαν κάποιοι έλεγχοι που αφορούσανε τα στοιχεία των χρηστών και συγκεκριμένα την αλλαγή κωδικού από το μενού αλλαγής στοιχείων το οποίο οδηγούσε στην δυσλειτουργία του συστήματος. Έχει γίνει κατάλληλη αναφορά. Η διόρθωση του σφάλματος κόστισε 1 εργατοώρα.</w:t></w:r></w:p>'
$textP.Range.InsertXML($xmlText)

# Step 3: add a fresh trailing empty paragraph (same style, no text) after that.
$textP = $d.Paragraphs.Item($count)
$textP.Range.InsertParagraphAfter()
$trailP = $d.Paragraphs.Item($count + 1)
$xmlTrail = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a3"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
$trailP.Range.InsertXML($xmlTrail)

# Step 4: the document still needs exactly one terminal paragraph mark, so fold
# the old bare paragraph away - the new trailing paragraph from step 3 survives.
$trailP = $d.Paragraphs.Item($count + 1)
$bareP = $d.Paragraphs.Item($count + 2)
$delRange = $d.Range($trailP.Range.End - 1, $bareP.Range.End)
$delRange.Delete()
